# Applies the "added extent report in context with threading" edit:
# adds a new test-data participant "arnav" (AddCustomerTest block) and
# "arnav k" (OpenAccountTest block) to the TestData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Insert a new row above row 5 - shifts the OpenAccountTest block
# (previously rows 6-9) down to rows 7-10, leaving row 5 free for the
# new AddCustomerTest data row while row 6 stays blank as a separator.
$ws.Rows("5:5").Insert()

# Append the new OpenAccountTest row (row 11, right after the existing
# data which now ends at row 10). Written first so "arnav k" is added
# to the shared-string table ahead of "arnav".
$ws.Range("A11").Value = "Y"
$ws.Range("B11").Value = "arnav k"
$ws.Range("C11").Value = "Dollar"
$ws.Range("D11").Value = "firefox"

# Fill in the new AddCustomerTest row (row 5)
$ws.Range("A5").Value = "Y"
$ws.Range("B5").Value = "arnav"
$ws.Range("C5").Value = "k"
$ws.Range("D5").Value = "X7878"
$ws.Range("E5").Value = "firefox"

# Update the selection to match the author's final cursor position
$ws.Range("C8").Select()
